$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3415
$ws.Range("J43").Value = 3582.8
$ws.Range("L43").Value = 3582.8
$ws.Range("N43").Value = -3720.8

$ws.Range("H76").Value = 9922.429
$ws.Range("I76").Value = 14670
$ws.Range("K76").Value = 14670
$ws.Range("M76").Value = -14355

$ws.Range("H79").Value = 9922.429
$ws.Range("I79").Value = 14670
$ws.Range("K79").Value = 14670
$ws.Range("M79").Value = -13578

$ws.Range("H92").Value = 3407
$ws.Range("I92").Value = 1169.1428
$ws.Range("J92").Value = 4712.4165
$ws.Range("K92").Value = 1169.1428
$ws.Range("L92").Value = 4712.4165
$ws.Range("M92").Value = 78.85719999999992
$ws.Range("N92").Value = -7208.4165

$ws.Range("H94").Value = 6331.3335
$ws.Range("I94").Value = 6331.3335
$ws.Range("K94").Value = 6331.3335
$ws.Range("M94").Value = -5880.3335

$ws.Range("H98").Value = 1738.75
$ws.Range("I98").Value = 1470.1724
$ws.Range("K98").Value = 1470.1724
$ws.Range("M98").Value = 27.82760000000007

$ws.Range("H99").Value = 3812.1428
$ws.Range("J99").Value = 5259.8
$ws.Range("L99").Value = 15779.4
$ws.Range("N99").Value = -18775.4

$ws.Range("H106").Value = 9104.3125
$ws.Range("I106").Value = 8984.929
$ws.Range("K106").Value = 8984.929
$ws.Range("M106").Value = -8353.929

$ws.Range("H107").Value = 2318.3
$ws.Range("J107").Value = 4921.5
$ws.Range("L107").Value = 4921.5
$ws.Range("N107").Value = -8761.5

$ws.Range("H122").Value = 1738.75
$ws.Range("I122").Value = 1470.1724
$ws.Range("K122").Value = 4410.5172
$ws.Range("M122").Value = -1960.5172

$ws.Range("H125").Value = 8115.2856
$ws.Range("I125").Value = 3850
$ws.Range("K125").Value = 34650
$ws.Range("M125").Value = -32190

$ws.Range("H135").Value = 2375.6785
$ws.Range("I135").Value = 753.087
$ws.Range("K135").Value = 6777.782999999999
$ws.Range("M135").Value = -4242.782999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2137.4285
$ws.Range("J2").Value = 3914.5
$ws.Range("L2").Value = 3914.5
$ws.Range("N2").Value = -4140.5

$ws.Range("H5").Value = 307.27274
$ws.Range("I5").Value = 326
$ws.Range("K5").Value = 326
$ws.Range("M5").Value = -214

$ws.Range("H32").Value = 4089.6553
$ws.Range("J32").Value = 1013
$ws.Range("L32").Value = 1013
$ws.Range("N32").Value = -1587

$ws.Range("H102").Value = 4622.6665
$ws.Range("I102").Value = 3444
$ws.Range("K102").Value = 3444
$ws.Range("M102").Value = -1822

$ws.Range("H116").Value = 2137.4285
$ws.Range("J116").Value = 3914.5
$ws.Range("L116").Value = 3914.5
$ws.Range("N116").Value = -8502.5

$ws.Range("H122").Value = 3159.9443
$ws.Range("I122").Value = 2645.3076
$ws.Range("K122").Value = 7935.9228
$ws.Range("M122").Value = -5485.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2137.4285
$ws.Range("J3").Value = 3914.5
$ws.Range("L3").Value = 3914.5
$ws.Range("N3").Value = -4142.5

$ws.Range("H4").Value = 307.27274
$ws.Range("I4").Value = 326
$ws.Range("K4").Value = 326
$ws.Range("M4").Value = -211

$ws.Range("H20").Value = 1523.2593
$ws.Range("I20").Value = 1273.7778
$ws.Range("J20").Value = 2022.2222
$ws.Range("K20").Value = 1273.7778
$ws.Range("L20").Value = 2022.2222
$ws.Range("M20").Value = -1026.7778
$ws.Range("N20").Value = -2516.2222

$ws.Range("H80").Value = 1018.05
$ws.Range("I80").Value = 442.5
$ws.Range("K80").Value = 442.5
$ws.Range("M80").Value = 555.5

$ws.Range("H83").Value = 1018.05
$ws.Range("I83").Value = 442.5
$ws.Range("K83").Value = 2212.5
$ws.Range("M83").Value = 2779.5

$ws.Range("H86").Value = 2439.8
$ws.Range("I86").Value = 1435.7368
$ws.Range("K86").Value = 1435.7368
$ws.Range("M86").Value = -312.7367999999999

$ws.Range("H89").Value = 2439.8
$ws.Range("I89").Value = 1435.7368
$ws.Range("K89").Value = 7178.683999999999
$ws.Range("M89").Value = -1562.683999999999

$ws.Range("H105").Value = 765716.6
$ws.Range("I105").Value = 1205889.5
$ws.Range("K105").Value = 1205889.5
$ws.Range("M105").Value = -1204142.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 43.57143
$ws.Range("I7").Value = 43.57143
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 43.57143
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 69.42857000000001
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 1351.1111
$ws.Range("I22").Value = 678.25
$ws.Range("J22").Value = 1889.4
$ws.Range("K22").Value = 678.25
$ws.Range("L22").Value = 1889.4
$ws.Range("M22").Value = -328.25
$ws.Range("N22").Value = -2589.4

$ws.Range("H31").Value = 31568616
$ws.Range("I31").Value = 43481020
$ws.Range("J31").Value = 1125801.6
$ws.Range("K31").Value = 43481020
$ws.Range("L31").Value = 1125801.6
$ws.Range("M31").Value = -43480725
$ws.Range("N31").Value = -1126391.6

$ws.Range("H34").Value = 31568616
$ws.Range("I34").Value = 43481020
$ws.Range("J34").Value = 1125801.6
$ws.Range("K34").Value = 43481020
$ws.Range("L34").Value = 1125801.6
$ws.Range("M34").Value = -43480818
$ws.Range("N34").Value = -1126205.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 13046.2
$ws.Range("J130").Value = 14331.714
$ws.Range("L130").Value = 42995.142
$ws.Range("N130").Value = -53035.142

$ws.Range("H140").Value = 5900.143
$ws.Range("I140").Value = 2781.889
$ws.Range("K140").Value = 8345.667000000001
$ws.Range("M140").Value = -3165.667000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 50021896
$ws.Range("I21").Value = 100014790
$ws.Range("J21").Value = 28999
$ws.Range("K21").Value = 100014790
$ws.Range("L21").Value = 28999
$ws.Range("M21").Value = -100014617
$ws.Range("N21").Value = -29345

$ws.Range("H30").Value = 50021896
$ws.Range("I30").Value = 100014790
$ws.Range("J30").Value = 28999
$ws.Range("K30").Value = 100014790
$ws.Range("L30").Value = 28999
$ws.Range("M30").Value = -100014685
$ws.Range("N30").Value = -29209

$ws.Range("H97").Value = 521.0769
$ws.Range("J97").Value = 554.1429000000001
$ws.Range("L97").Value = 554.1429000000001
$ws.Range("N97").Value = -1546.1429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8112.154
$ws.Range("I7").Value = 8059.8184
$ws.Range("K7").Value = 8059.8184
$ws.Range("M7").Value = -7947.8184

$ws.Range("H122").Value = 3561.6274
$ws.Range("I122").Value = 3434.9788
$ws.Range("K122").Value = 10304.9364
$ws.Range("M122").Value = -7854.936399999999

$ws.Range("H126").Value = 8112.154
$ws.Range("I126").Value = 8059.8184
$ws.Range("K126").Value = 24179.4552
$ws.Range("M126").Value = -21709.4552

$ws.Range("H132").Value = 5533.846
$ws.Range("I132").Value = 3189.5
$ws.Range("K132").Value = 9568.5
$ws.Range("M132").Value = -7038.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1524.2
$ws.Range("J23").Value = 5222
$ws.Range("L23").Value = 5222
$ws.Range("N23").Value = -5680

$ws.Range("H31").Value = 8000
$ws.Range("I31").Value = 8000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 8000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -7652
$ws.Range("N31").ClearContents()

$ws.Range("H37").Value = 35832.5
$ws.Range("I37").Value = 27000
$ws.Range("J37").Value = 79995
$ws.Range("K37").Value = 27000
$ws.Range("L37").Value = 79995
$ws.Range("M37").Value = -26797
$ws.Range("N37").Value = -80401

$ws.Range("H122").Value = 3483.2222
$ws.Range("I122").Value = 3070.125
$ws.Range("K122").Value = 9210.375
$ws.Range("M122").Value = -6760.375

$ws.Range("H126").Value = 8622.25
$ws.Range("I126").Value = 8622.25
$ws.Range("K126").Value = 25866.75
$ws.Range("M126").Value = -23396.75
